$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J so the old J column (and its data/width) shifts to K,
# then give the new J1 header its own label ("地區別代號").
$ws.Columns("J:J").Insert()
$ws.Range("J1").Value = "地區別代號"

# Re-apply the (approximate) column widths for the two affected columns.
$ws.Columns("J:J").ColumnWidth = 12.36328125
$ws.Columns("K:K").ColumnWidth = 10.453125

# Update the workbook-level Database defined name to include the new column.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Database") {
        $n.RefersTo = "=LAW7N1Ppp!`$A`$1:`$K`$1"
    }
}

# Move the active selection to match the author's final cursor position.
$ws.Range("J5").Select() | Out-Null
